# Apply updated TPM values to the Fgf16-Fgfr2 LR-pair sheet.
# - Rename the "Inflammatory-Mac" cluster label to "Resolving-Mac"
#   (Excel automatically updates the shared-string table for every cell
#   that referenced the old text when the new text is assigned).
# - Refresh the recomputed NATMI statistics (columns G-T) for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.1628603333333333  # G2
$ws.Cells.Item(2, 8).Value = 0.488581  # H2
$ws.Cells.Item(2, 9).Value = 0.1212602171476209  # I2
$ws.Cells.Item(2, 10).Value = 0.171492704432097  # J2
$ws.Cells.Item(2, 11).Value = 3  # K2
$ws.Cells.Item(2, 12).Value = 1  # L2
$ws.Cells.Item(2, 13).Value = 0.2858606666666667  # M2
$ws.Cells.Item(2, 14).Value = 0.857582  # N2
$ws.Cells.Item(2, 15).Value = 0.0687156860066334  # O2
$ws.Cells.Item(2, 16).Value = 0.06932858672617494  # P2
$ws.Cells.Item(2, 17).Value = 0.04655536346022222  # Q2
$ws.Cells.Item(2, 18).Value = 0.4189982711419999  # R2
$ws.Cells.Item(2, 19).Value = 0.008332479006612104  # S2
$ws.Cells.Item(2, 20).Value = 0.01188934683212692  # T2
# Row 3
$ws.Cells.Item(3, 7).Value = 0.1628603333333333  # G3
$ws.Cells.Item(3, 8).Value = 0.488581  # H3
$ws.Cells.Item(3, 9).Value = 0.1212602171476209  # I3
$ws.Cells.Item(3, 10).Value = 0.171492704432097  # J3
$ws.Cells.Item(3, 15).Value = 0.90464312565499  # O3
$ws.Cells.Item(3, 16).Value = 0.9127119736118995  # P3
$ws.Cells.Item(3, 17).Value = 0.6129021183401111  # Q3
$ws.Cells.Item(3, 18).Value = 5.516119065061  # R3
$ws.Cells.Item(3, 19).Value = 0.1096972218580266  # S3
$ws.Cells.Item(3, 20).Value = 0.1565234447222614  # T3
# Row 4
$ws.Cells.Item(4, 4).Value = "MuSCs"  # D4
$ws.Cells.Item(4, 7).Value = 0.1628603333333333  # G4
$ws.Cells.Item(4, 8).Value = 0.488581  # H4
$ws.Cells.Item(4, 9).Value = 0.1212602171476209  # I4
$ws.Cells.Item(4, 10).Value = 0.171492704432097  # J4
$ws.Cells.Item(4, 11).Value = 2  # K4
$ws.Cells.Item(4, 12).Value = 1  # L4
$ws.Cells.Item(4, 13).Value = 0.110331  # M4
$ws.Cells.Item(4, 14).Value = 0.220662  # N4
$ws.Cells.Item(4, 15).Value = 0.02652155835639462  # O4
$ws.Cells.Item(4, 16).Value = 0.01783874265571248  # P4
$ws.Cells.Item(4, 17).Value = 0.017968543437  # Q4
$ws.Cells.Item(4, 18).Value = 0.107811260622  # R4
$ws.Cells.Item(4, 19).Value = 0.003216009925389712  # S4
$ws.Cells.Item(4, 20).Value = 0.00305921422169634  # T4
# Row 5
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"  # D5
$ws.Cells.Item(5, 7).Value = 0.1628603333333333  # G5
$ws.Cells.Item(5, 8).Value = 0.488581  # H5
$ws.Cells.Item(5, 9).Value = 0.1212602171476209  # I5
$ws.Cells.Item(5, 10).Value = 0.171492704432097  # J5
$ws.Cells.Item(5, 11).Value = 1  # K5
$ws.Cells.Item(5, 12).Value = 0.3333333333333333  # L5
$ws.Cells.Item(5, 13).Value = 0.0004976666666666667  # M5
$ws.Cells.Item(5, 14).Value = 0.001493  # N5
$ws.Cells.Item(5, 15).Value = 0.0001196299819817856  # O5
$ws.Cells.Item(5, 16).Value = 0.0001206970062130259  # P5
$ws.Cells.Item(5, 17).Value = 0.00008105015922222222  # Q5
$ws.Cells.Item(5, 18).Value = 0.000729451433  # R5
$ws.Cells.Item(5, 19).Value = 0.0000145063575924773  # S5
$ws.Cells.Item(5, 20).Value = 0.00002069865601232943  # T5
# Row 6
$ws.Cells.Item(6, 7).Value = 1.1802045  # G6
$ws.Cells.Item(6, 8).Value = 2.360409  # H6
$ws.Cells.Item(6, 9).Value = 0.878739782852379  # I6
$ws.Cells.Item(6, 10).Value = 0.828507295567903  # J6
$ws.Cells.Item(6, 11).Value = 3  # K6
$ws.Cells.Item(6, 12).Value = 1  # L6
$ws.Cells.Item(6, 13).Value = 0.2858606666666667  # M6
$ws.Cells.Item(6, 14).Value = 0.857582  # N6
$ws.Cells.Item(6, 15).Value = 0.0687156860066334  # O6
$ws.Cells.Item(6, 16).Value = 0.06932858672617494  # P6
$ws.Cells.Item(6, 17).Value = 0.337374045173  # Q6
$ws.Cells.Item(6, 18).Value = 2.024244271038  # R6
$ws.Cells.Item(6, 19).Value = 0.06038320700002129  # S6
$ws.Cells.Item(6, 20).Value = 0.05743923989404802  # T6
# Row 7
$ws.Cells.Item(7, 7).Value = 1.1802045  # G7
$ws.Cells.Item(7, 8).Value = 2.360409  # H7
$ws.Cells.Item(7, 9).Value = 0.878739782852379  # I7
$ws.Cells.Item(7, 10).Value = 0.828507295567903  # J7
$ws.Cells.Item(7, 15).Value = 0.90464312565499  # O7
$ws.Cells.Item(7, 16).Value = 0.9127119736118995  # P7
$ws.Cells.Item(7, 17).Value = 4.4415348005215  # Q7
$ws.Cells.Item(7, 18).Value = 26.649208803129  # R7
$ws.Cells.Item(7, 19).Value = 0.7949459037969634  # S7
$ws.Cells.Item(7, 20).Value = 0.7561885288896382  # T7
# Row 8
$ws.Cells.Item(8, 4).Value = "MuSCs"  # D8
$ws.Cells.Item(8, 7).Value = 1.1802045  # G8
$ws.Cells.Item(8, 8).Value = 2.360409  # H8
$ws.Cells.Item(8, 9).Value = 0.878739782852379  # I8
$ws.Cells.Item(8, 10).Value = 0.828507295567903  # J8
$ws.Cells.Item(8, 11).Value = 2  # K8
$ws.Cells.Item(8, 12).Value = 1  # L8
$ws.Cells.Item(8, 13).Value = 0.110331  # M8
$ws.Cells.Item(8, 14).Value = 0.220662  # N8
$ws.Cells.Item(8, 15).Value = 0.02652155835639462  # O8
$ws.Cells.Item(8, 16).Value = 0.01783874265571248  # P8
$ws.Cells.Item(8, 17).Value = 0.1302131426895  # Q8
$ws.Cells.Item(8, 18).Value = 0.520852570758  # R8
$ws.Cells.Item(8, 19).Value = 0.0233055484310049  # S8
$ws.Cells.Item(8, 20).Value = 0.01477952843401614  # T8
# Row 9
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"  # D9
$ws.Cells.Item(9, 7).Value = 1.1802045  # G9
$ws.Cells.Item(9, 8).Value = 2.360409  # H9
$ws.Cells.Item(9, 9).Value = 0.878739782852379  # I9
$ws.Cells.Item(9, 10).Value = 0.828507295567903  # J9
$ws.Cells.Item(9, 11).Value = 1  # K9
$ws.Cells.Item(9, 12).Value = 0.3333333333333333  # L9
$ws.Cells.Item(9, 13).Value = 0.0004976666666666667  # M9
$ws.Cells.Item(9, 14).Value = 0.001493  # N9
$ws.Cells.Item(9, 15).Value = 0.0001196299819817856  # O9
$ws.Cells.Item(9, 16).Value = 0.0001206970062130259  # P9
$ws.Cells.Item(9, 17).Value = 0.0005873484394999999  # Q9
$ws.Cells.Item(9, 18).Value = 0.003524090637  # R9
$ws.Cells.Item(9, 19).Value = 0.0001051236243893083  # S9
$ws.Cells.Item(9, 20).Value = 0.00009999835020069649  # T9
